# Updates the cryptos price/volume snapshot (Sheet1, rows 2-51).
# Columns: B=Coin, C=Link, D=Price, E=Volume(1h)
# Values that look numeric to Excel (single-dot decimals such as "1.008")
# are written with a leading apostrophe so they stay text, matching the
# original inlineStr cell type instead of being auto-converted to a number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.215.93"
$ws.Range("E2").Value = "  +0.46%  "
$ws.Range("E3").Value = "  +0.13%  "
$ws.Range("D4").Value = "'1.008"
$ws.Range("E4").Value = "  +0.54%  "
$ws.Range("D5").Value = "'218.03"
$ws.Range("E5").Value = "  -0.09%  "
$ws.Range("D6").Value = "'0.5310"
$ws.Range("E6").Value = "  -0.04%  "
$ws.Range("D7").Value = "'1.008"
$ws.Range("E7").Value = "  +0.47%  "
$ws.Range("D8").Value = "'0.2625"
$ws.Range("E8").Value = "  +0.20%  "
$ws.Range("D9").Value = "'0.06340"
$ws.Range("E9").Value = "  +0.65%  "
$ws.Range("E10").Value = "  +0.38%  "
$ws.Range("D11").Value = "'0.07834"
$ws.Range("E11").Value = "  +1.10%  "
$ws.Range("D12").Value = "'4.526"
$ws.Range("E12").Value = "  +0.98%  "
$ws.Range("D13").Value = "1.633.92"
$ws.Range("E13").Value = "  -1.37%  "
$ws.Range("D14").Value = "1.883.78"
$ws.Range("E14").Value = "  +0.17%  "
$ws.Range("D15").Value = "'0.5503"
$ws.Range("E15").Value = "  +0.85%  "
$ws.Range("D16").Value = "0.0₅8166"
$ws.Range("E16").Value = "  +0.54%  "
$ws.Range("D17").Value = "'65.42"
$ws.Range("E17").Value = "  +0.22%  "
$ws.Range("D18").Value = "26.183.58"
$ws.Range("E18").Value = "  +0.26%  "
$ws.Range("D19").Value = "'1.008"
$ws.Range("E19").Value = "  +0.53%  "
$ws.Range("D20").Value = "'4.608"
$ws.Range("E20").Value = "  +0.95%  "
$ws.Range("D21").Value = "'191.43"
$ws.Range("D22").Value = "'10.10"
$ws.Range("E22").Value = "  +0.50%  "
$ws.Range("D23").Value = "'6.022"
$ws.Range("E23").Value = "  +0.16%  "
$ws.Range("D24").Value = "'1.009"
$ws.Range("E24").Value = "  +0.49%  "
$ws.Range("D25").Value = "'143.55"
$ws.Range("E25").Value = "  +2.71%  "
$ws.Range("D26").Value = "'0.1220"
$ws.Range("E26").Value = "  -2.15%  "
$ws.Range("D27").Value = "'7.210"
$ws.Range("E27").Value = "  -0.70%  "
$ws.Range("E28").Value = "  -1.26%  "
$ws.Range("D29").Value = "'1.472"
$ws.Range("E29").Value = "  +2.76%  "
$ws.Range("D30").Value = "'0.05780"
$ws.Range("E30").Value = "  -2.71%  "
$ws.Range("D31").Value = "'1.276"
$ws.Range("E31").Value = "  -0.19%  "
$ws.Range("D32").Value = "'3.558"
$ws.Range("E32").Value = "  +1.53%  "
$ws.Range("D33").Value = "'3.273"
$ws.Range("E33").Value = "  +0.75%  "
$ws.Range("D34").Value = "'1.597"
$ws.Range("E34").Value = "  +3.33%  "
$ws.Range("D35").Value = "'2.816"
$ws.Range("E35").Value = "  +2.07%  "
$ws.Range("D36").Value = "'0.9517"
$ws.Range("E36").Value = "  +1.04%  "
$ws.Range("D37").Value = "'2.427"
$ws.Range("E37").Value = "  +0.64%  "
$ws.Range("D38").Value = "'0.5770"
$ws.Range("E38").Value = "  +1.88%  "
$ws.Range("D39").Value = "'0.01601"
$ws.Range("E39").Value = "  -0.46%  "
$ws.Range("E40").Value = "  -0.56%  "
$ws.Range("E41").Value = "  +0.49%  "
$ws.Range("B42").Value = "PaxDollar"
$ws.Range("C42").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D42").Value = "'1.007"
$ws.Range("E42").Value = "  +0.48%  "
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "1.045.03"
$ws.Range("E43").Value = "  +3.45%  "
$ws.Range("D44").Value = "'103.95"
$ws.Range("E44").Value = "  +2.94%  "
$ws.Range("D45").Value = "1.796.70"
$ws.Range("E45").Value = "  +0.02%  "
$ws.Range("D46").Value = "'56.87"
$ws.Range("E46").Value = "  -0.34%  "
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").Value = "0.0₈105"
$ws.Range("E47").Value = "  -0.61%  "
$ws.Range("B48").Value = "Frax"
$ws.Range("C48").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D48").Value = "'1.007"
$ws.Range("E48").Value = "  +0.65%  "
$ws.Range("E49").Value = "  +1.69%  "
$ws.Range("D50").Value = "'7.876"
$ws.Range("E50").Value = "  -0.08%  "
$ws.Range("D51").Value = "'0.05156"
$ws.Range("E51").Value = "  +0.09%  "
